# Auto-applies the diff described for 杭州-漫展信息.xlsx (update at 456a3b4)
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")

# Insert new row 37 for the new "ACG CLUB" event, shifting rows 37-45 down to 38-46
$ws1.Rows.Item(37).Insert()

# Copy formatting (border/alignment/font) from A36 onto the newly inserted A37 cell
$ws1.Range("A36").Copy()
$ws1.Range("A37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Range("A37").Value = 36
$b37 = $ws1.Range("B37")
$b37.NumberFormat = "@"
$b37.Value = "2024-07-27"
$ws1.Range("C37").Value = "杭州·ACG CLUB动漫游戏嘉年华"
$ws1.Range("D37").Value = "中心路1号 白蓝地文创街区"
$ws1.Range("E37").Value = "2024.07.27 10:00-07.27 17:00"
$ws1.Range("F37").Value = 2
$ws1.Range("G37").Value = 68.88
$ws1.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=86265"
$ws1.Range("I37").Value = "//i1.hdslb.com/bfs/openplatform/202405/XBRfeQwu1716533419093.jpeg"

# Refresh "want-to-go" counts (column F) for unaffected rows
$ws1.Range("F3").Value = 251
$ws1.Range("F4").Value = 257
$ws1.Range("F5").Value = 1795
$ws1.Range("F7").Value = 530
$ws1.Range("F8").Value = 4870
$ws1.Range("F13").Value = 331
$ws1.Range("F14").Value = 1293
$ws1.Range("F16").Value = 1887
$ws1.Range("F17").Value = 3001
$ws1.Range("F18").Value = 1832
$ws1.Range("F22").Value = 65
$ws1.Range("F25").Value = 313
$ws1.Range("F26").Value = 34
$ws1.Range("F27").Value = 3026
$ws1.Range("F28").Value = 1034
$ws1.Range("F29").Value = 2516
$ws1.Range("F31").Value = 1364
$ws1.Range("F32").Value = 3659
$ws1.Range("F33").Value = 95
$ws1.Range("F34").Value = 897
$ws1.Range("F35").Value = 434
$ws1.Range("F36").Value = 1153

# Refresh "want-to-go" counts (column F) for rows shifted down by the insertion above
$ws1.Range("F38").Value = 938
$ws1.Range("F39").Value = 1200
$ws1.Range("F41").Value = 893
$ws1.Range("F42").Value = 574
$ws1.Range("F43").Value = 306
$ws1.Range("F45").Value = 292
$ws1.Range("F46").Value = 3501

# ---- Sheet 2: 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 17
$ws2.Range("F5").Value = 2

# ---- Sheet 4: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 251
$ws4.Range("F4").Value = 257
$ws4.Range("F6").Value = 1795
$ws4.Range("F8").Value = 530
$ws4.Range("F9").Value = 4870
$ws4.Range("F11").Value = 17
$ws4.Range("F13").Value = 331
$ws4.Range("F14").Value = 1293
$ws4.Range("F15").Value = 3001
$ws4.Range("F17").Value = 1832
$ws4.Range("F25").Value = 65
$ws4.Range("F27").Value = 313
$ws4.Range("F28").Value = 3026
$ws4.Range("F30").Value = 1034
$ws4.Range("F31").Value = 2516
$ws4.Range("F32").Value = 1364
$ws4.Range("F33").Value = 3659
$ws4.Range("F35").Value = 95
$ws4.Range("F36").Value = 897
$ws4.Range("F37").Value = 1153
$ws4.Range("F38").Value = 938
$ws4.Range("F40").Value = 1200
$ws4.Range("F41").Value = 893
$ws4.Range("F42").Value = 574
$ws4.Range("F47").Value = 292
$ws4.Range("F48").Value = 3501

